$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("D1").Comment
$txt = $c.Text()
$c.Text($txt)
Write-Output "done"
